$wb = $excel.ActiveWorkbook

# Add a new worksheet "EPFImporter" at the end of the workbook (after the
# last existing sheet), which also makes it the active sheet / active tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "EPFImporter"

# Fill in the sketch of the EPFImporter / EPFIngester / EPFParser / EPFDbUtil
# interfaces, in the same order the values were originally entered (this
# keeps the shared-string table append order identical to the source file).
$ws.Range("A1").Value = "EPFImporter Interface Definitions"
$ws.Range("A3").Value = "EPFImporter"
$ws.Range("B4").Value = "Description"
$ws.Range("B5").Value = "Load Configuration Files"
$ws.Range("B6").Value = "Load Command Line Parameters"
$ws.Range("C7").Value = "Import Directory"
$ws.Range("C8").Value = "Import Files"
$ws.Range("C9").Value = "Import Whitelist"
$ws.Range("C10").Value = "Import Blacklist"
$ws.Range("B13").Value = "For each file to import"
$ws.Range("C11").Value = "Load Thread Pool Maximum"
$ws.Range("A16").Value = "EPFIngester"
$ws.Range("C14").Value = "ExecutionQueue.add(new EPFIngester(importFile))"
$ws.Range("B17").Value = "new EPFParser(new EPFFileReader(importFile))"
$ws.Range("B19").Value = "parseTableName"
$ws.Range("B20").Value = "parseColumnsAndTypes"
$ws.Range("B21").Value = "parsePrimaryKey"
$ws.Range("B22").Value = "seekRecord()"
$ws.Range("B23").Value = "totalRecords()"
$ws.Range("B24").Value = "nextRecord()"
$ws.Range("A27").Value = "EPFDbUtil"
$ws.Range("B28").Value = "initTable(tableName,importType,totalRecords)"
$ws.Range("B29").Value = "createTable(tableName,columnsAndTypes,primaryKeys)"
$ws.Range("B30").Value = "insertRow(List<String> values)"
$ws.Range("B31").Value = "finalizeTable()"

# Match the recorded view state for the new sheet: scrolled so row 8 is the
# top-left visible row, with A33 selected (the cell just past the last used
# row, where the user would naturally continue typing).
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("A33").Select()
